$wb = $excel.ActiveWorkbook

$xlShiftToRight = -4161
$xlShiftToLeft  = -4159
$xlPasteFormats = -4122

# =========================================================================
# Sheet "Login" (sheet1): insert "S.No" / "Test Description" columns
# before the existing Uname/pwd columns.
# =========================================================================
$ws1 = $wb.Worksheets.Item("Login")

# Push existing A:B (Uname/pwd table) two columns to the right -> C:D,
# carrying their original formatting (border/fill/bold header) with them.
$ws1.Range("A1:B3").Insert($xlShiftToRight)

# Stamp the new A/B columns with the same look as the table they extend.
$ws1.Range("C1").Copy()
$ws1.Range("A1:B1").PasteSpecial($xlPasteFormats)
$ws1.Range("C2").Copy()
$ws1.Range("A2:B3").PasteSpecial($xlPasteFormats)
$ws1.Application.CutCopyMode = $false

$ws1.Range("A1").Value = "S.No"
$ws1.Range("B1").Value = "Test Description"
$ws1.Range("A2").Value = "TC_0001"
$ws1.Range("B2").Value = "Verify User is able to login"
$ws1.Range("A3").Value = "TC_0002"
$ws1.Range("B3").Value = "Verify User is able to login and logout"

$ws1.Columns.Item(1).ColumnWidth = 8.140625
$ws1.Columns.Item(2).ColumnWidth = 34.85546875
$ws1.Columns.Item(4).ColumnWidth = 6.85546875

$ws1.Range("C7").Select()

# =========================================================================
# Sheet "Create" (sheet2): drop the companyName/firstName/lastName columns,
# add "S.No" / description columns in front of Uname/pwd.
# =========================================================================
$ws2 = $wb.Worksheets.Item("Create")

$ws2.Columns("C:E").Delete()
$ws2.Range("A1:B3").Insert($xlShiftToRight)

$ws2.Range("C1").Copy()
$ws2.Range("A1:B1").PasteSpecial($xlPasteFormats)
$ws2.Range("C2").Copy()
$ws2.Range("A2:B3").PasteSpecial($xlPasteFormats)
$ws2.Application.CutCopyMode = $false

$ws2.Range("A1").Value = "S.No"
$ws2.Range("B1").Value = "Uname"
$ws2.Range("A2").Value = "TC_0001"
$ws2.Range("B2").Value = "Verify User is able to create lead using DemoSalesManager login"
$ws2.Range("A3").Value = "TC_0002"
$ws2.Range("B3").Value = "Verify User is able to create lead using DemoCSR login"

$ws2.Columns.Item(1).ColumnWidth = 8.140625
$ws2.Columns.Item(2).ColumnWidth = 59.28515625
$ws2.Columns.Item(4).ColumnWidth = 6.85546875

$ws2.Range("B13").Select()

# =========================================================================
# Sheet "Find" (sheet3): same shape as "Create" with different wording.
# =========================================================================
$ws3 = $wb.Worksheets.Item("Find")

$ws3.Columns("C:E").Delete()
$ws3.Range("A1:B3").Insert($xlShiftToRight)

$ws3.Range("C1").Copy()
$ws3.Range("A1:B1").PasteSpecial($xlPasteFormats)
$ws3.Range("C2").Copy()
$ws3.Range("A2:B3").PasteSpecial($xlPasteFormats)
$ws3.Application.CutCopyMode = $false

$ws3.Range("A1").Value = "S.No"
$ws3.Range("B1").Value = "Uname"
$ws3.Range("A2").Value = "TC_0001"
$ws3.Range("B2").Value = "Verify User is able to create and find lead using DemoSalesManager login"
$ws3.Range("A3").Value = "TC_0002"
$ws3.Range("B3").Value = "Verify User is able to create and find lead using DemoCSR login"

$ws3.Columns.Item(2).ColumnWidth = 67.28515625

$ws3.Range("B7").Select()
$ws3.Activate()
